$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4035.6667
$ws.Range("I76").Value = 4035.6667
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 4035.6667
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -3720.6667
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 4035.6667
$ws.Range("I79").Value = 4035.6667
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 4035.6667
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -2943.6667
$ws.Range("N79").ClearContents()

$ws.Range("H82").Value = 7486
$ws.Range("I82").Value = 7486
$ws.Range("K82").Value = 22458
$ws.Range("M82").Value = -22052

$ws.Range("H85").Value = 7486
$ws.Range("I85").Value = 7486
$ws.Range("K85").Value = 22458
$ws.Range("M85").Value = -21054

$ws.Range("H88").Value = 25561.117
$ws.Range("J88").Value = 47101.668
$ws.Range("L88").Value = 47101.668
$ws.Range("N88").Value = -47913.668

$ws.Range("H91").Value = 25561.117
$ws.Range("J91").Value = 47101.668
$ws.Range("L91").Value = 47101.668
$ws.Range("N91").Value = -49909.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1321.8667
$ws.Range("I88").Value = 1019.5714
$ws.Range("J88").Value = 1586.375
$ws.Range("K88").Value = 1019.5714
$ws.Range("L88").Value = 1586.375
$ws.Range("M88").Value = -613.5714
$ws.Range("N88").Value = -2398.375

$ws.Range("H91").Value = 1321.8667
$ws.Range("I91").Value = 1019.5714
$ws.Range("J91").Value = 1586.375
$ws.Range("K91").Value = 1019.5714
$ws.Range("L91").Value = 1586.375
$ws.Range("M91").Value = 384.4286
$ws.Range("N91").Value = -4394.375

$ws.Range("H96").Value = 41210.95
$ws.Range("J96").Value = 41210.95
$ws.Range("L96").Value = 41210.95
$ws.Range("N96").Value = -46702.95

$ws.Range("H122").Value = 2721.0938
$ws.Range("I122").Value = 2434.1155
$ws.Range("J122").Value = 3964.6667
$ws.Range("K122").Value = 7302.3465
$ws.Range("L122").Value = 11894.0001
$ws.Range("M122").Value = -4852.3465
$ws.Range("N122").Value = -16794.0001

$ws.Range("H132").Value = 7840.1763
$ws.Range("I132").Value = 9511.944
$ws.Range("K132").Value = 28535.832
$ws.Range("M132").Value = -26005.832

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2446.6553
$ws.Range("I86").Value = 2419.923
$ws.Range("K86").Value = 2419.923
$ws.Range("M86").Value = -1296.923

$ws.Range("H89").Value = 2446.6553
$ws.Range("I89").Value = 2419.923
$ws.Range("K89").Value = 12099.614999999998
$ws.Range("M89").Value = -6483.614999999998

$ws.Range("H134").Value = 3080.3333
$ws.Range("I134").Value = 2568.64
$ws.Range("J134").Value = 4243.273
$ws.Range("K134").Value = 7705.92
$ws.Range("L134").Value = 12729.819
$ws.Range("M134").Value = -5170.92
$ws.Range("N134").Value = -17799.819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 6801.3
$ws.Range("I62").Value = 9948
$ws.Range("J62").Value = 6014.625
$ws.Range("K62").Value = 9948
$ws.Range("L62").Value = 6014.625
$ws.Range("M62").Value = -9324
$ws.Range("N62").Value = -7262.625

$ws.Range("H65").Value = 6801.3
$ws.Range("I65").Value = 9948
$ws.Range("J65").Value = 6014.625
$ws.Range("K65").Value = 49740
$ws.Range("L65").Value = 30073.125
$ws.Range("M65").Value = -46620
$ws.Range("N65").Value = -36313.125

$ws.Range("H98").Value = 66666
$ws.Range("J98").Value = 66666
$ws.Range("L98").Value = 66666
$ws.Range("N98").Value = -71158

$ws.Range("H100").Value = 73501.164
$ws.Range("I100").Value = 33709
$ws.Range("J100").Value = 81459.60000000001
$ws.Range("K100").Value = 33709
$ws.Range("L100").Value = 81459.60000000001
$ws.Range("M100").Value = -32627
$ws.Range("N100").Value = -83623.60000000001

$ws.Range("H107").Value = 2476.6052
$ws.Range("I107").Value = 1364.52
$ws.Range("K107").Value = 1364.52
$ws.Range("M107").Value = 555.48

$ws.Range("H124").Value = 45000
$ws.Range("J124").Value = 45000
$ws.Range("L124").Value = 45000
$ws.Range("N124").Value = -49910

$ws.Range("H132").Value = 6834.41
$ws.Range("I132").Value = 3926.0312
$ws.Range("J132").Value = 20129.857
$ws.Range("K132").Value = 11778.0936
$ws.Range("L132").Value = 60389.571
$ws.Range("M132").Value = -9248.0936
$ws.Range("N132").Value = -65449.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1351.7428
$ws.Range("J107").Value = 1368.8529
$ws.Range("L107").Value = 4106.5587
$ws.Range("N107").Value = -7946.5587

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5684.25
$ws.Range("I70").Value = 5912.5
$ws.Range("K70").Value = 5912.5
$ws.Range("M70").Value = -5642.5

$ws.Range("H73").Value = 5684.25
$ws.Range("I73").Value = 5912.5
$ws.Range("K73").Value = 5912.5
$ws.Range("M73").Value = -4976.5

$ws.Range("H102").Value = 23823.492
$ws.Range("I102").Value = 27780.83
$ws.Range("K102").Value = 27780.83
$ws.Range("M102").Value = -26158.83

$ws.Range("H109").Value = 37996
$ws.Range("J109").Value = 37996
$ws.Range("L109").Value = 37996
$ws.Range("N109").Value = -40076

$ws.Range("H132").Value = 3274.3972
$ws.Range("I132").Value = 2831.8333
$ws.Range("J132").Value = 5317
$ws.Range("K132").Value = 8495.499899999999
$ws.Range("L132").Value = 15951
$ws.Range("M132").Value = -5965.499899999999
$ws.Range("N132").Value = -21011

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4013
$ws.Range("I61").Value = 3951.3125
$ws.Range("K61").Value = 3951.3125
$ws.Range("M61").Value = -3749.3125

$ws.Range("H82").Value = 29413124
$ws.Range("I82").Value = 1517.84
$ws.Range("J82").Value = 111112030
$ws.Range("K82").Value = 1517.84
$ws.Range("L82").Value = 111112030
$ws.Range("M82").Value = -1156.84
$ws.Range("N82").Value = -111112752

$ws.Range("H85").Value = 29413124
$ws.Range("I85").Value = 1517.84
$ws.Range("J85").Value = 111112030
$ws.Range("K85").Value = 1517.84
$ws.Range("L85").Value = 111112030
$ws.Range("M85").Value = -269.8399999999999
$ws.Range("N85").Value = -111114526

$ws.Range("H93").Value = 821.9048
$ws.Range("I93").Value = 737.4666999999999
$ws.Range("K93").Value = 737.4666999999999
$ws.Range("M93").Value = 510.5333000000001

$ws.Range("H113").Value = 4013
$ws.Range("I113").Value = 3951.3125
$ws.Range("K113").Value = 3951.3125
$ws.Range("M113").Value = -1781.3125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 23435.75
$ws.Range("J28").Value = 21995
$ws.Range("L28").Value = 21995
$ws.Range("N28").Value = -22691

$ws.Range("H132").Value = 5132.162
$ws.Range("I132").Value = 3689.8823
$ws.Range("J132").Value = 6358.1
$ws.Range("K132").Value = 11069.6469
$ws.Range("L132").Value = 19074.3
$ws.Range("M132").Value = -8539.6469
$ws.Range("N132").Value = -24134.3
